# "send kode otp to email and revisi 27 mei"
# Replace the Kode_OTP column with an Email column (values become mailto
# hyperlinks), and move the active selection from E7 to C8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the Kode_OTP header to Email -----------------------------
$ws.Range("D1").Value = "Email"

# --- Replace the OTP codes with the alumni's e-mail addresses --------
$ws.Range("D2").Value = "arindrakeysha@gmail.com"
$ws.Range("D3").Value = "cascanekeysha@gmail.com"
$ws.Range("D4").Value = "dekuw85@gmail.com"
$ws.Range("D5").Value = "dinarullailil26@gmail.com"

# --- Turn each e-mail address into a mailto: hyperlink ----------------
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:arindrakeysha@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:cascanekeysha@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:dekuw85@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D5"), "mailto:dinarullailil26@gmail.com")

# --- Move the saved selection (cosmetic "revisi 27 mei" tweak) -------
[void]$ws.Range("C8").Select()
